# Add "Name on Account" column to the Tokenized Bank Accounts template.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New header cell F1, using the same style as the other grey header cells (A1/C1/D1/E1).
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Name on Account"

# Give the new "Name on Account" column a wider column, like the other text columns.
$ws.Columns.Item(6).ColumnWidth = 19.92

# Document the new column the same way the other header cells are documented.
$ws.Range("F1").AddComment("The name on the bank account")

# Keep selection consistent with the new last column, matching the template's convention
# of leaving the cursor on the most-recently-added header cell.
$ws.Range("F1").Select()
